$wb = $excel.ActiveWorkbook

# Update the efficiency-improvement values for heating, cooling/ventilation,
# and appliances (rows 2, 3, 6) on the PPEIdtIL sheet from 0.1 to 0.02.
$ws = $wb.Worksheets.Item("PPEIdtIL")
$ws.Range("B2").Value = 0.02
$ws.Range("D2").Value = 0.02
$ws.Range("B3").Value = 0.02
$ws.Range("D3").Value = 0.02
$ws.Range("B6").Value = 0.02
$ws.Range("D6").Value = 0.02

# Switch the active sheet/tab to PPEIdtIL and restore its selection.
$ws.Activate()
$ws.Range("L3").Select()

# Update the About sheet's selection (it is no longer the active tab).
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("E35").Select()

# Re-activate PPEIdtIL so it ends up as the active/selected tab.
$ws.Activate()
